$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.798.28'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '3.121.63'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.87'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.47'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  -3.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.386'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.41%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '3.119.28'
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  -3.04%  '
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.89'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.32%  '
$ws.Range("D16").Value = '91.604.32'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '3.697.61'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '3.139.85'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("E19").Value = '  +0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '447.81'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000203'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.144'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +27.91%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.233'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("E32").Value = '  -10.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.177'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -1.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.25'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.11'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '492.41'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.441'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.88%  '
$ws.Range("E43").Value = '  -6.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '157.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.702'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("E51").Value = '  -3.31%  '
